$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: merge a run-spanning range of text (identified by unique start/end
# substrings) back into a single run, without altering its visible text.
# The underlying engine auto-normalises (merges) adjacent same-formatted runs
# whenever a Range.Text assignment actually changes the stored text, so we
# briefly perturb the text and then restore it.
# ---------------------------------------------------------------------------
function Merge-Range($startNeedle, $endNeedle) {
    $r1 = $d.Content
    $r1.Find.Execute($startNeedle, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $start = $r1.Start

    $r2 = $d.Content
    $r2.Find.Execute($endNeedle, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $end = $r2.End

    $rng = $d.Range($start, $end)
    $orig = $rng.Text
    $rng.Text = $orig + "%"
    $rng2 = $d.Range($start, $start + $orig.Length + 1)
    $rng2.Text = $orig
}

# ---------------------------------------------------------------------------
# Change 1: "eficiencia" -> "eficacia", landing in its own run (the run that
# held the whole sentence gets split into three runs around it, matching the
# bookmark that sits right after).
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("eficiencia", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Text = ""
$rng.InsertAfter("eficacia")
$newWord = $d.Range($rng.Start, $rng.Start + 8)
$newWord.Bold = $true
$newWord.Bold = $false

# The "_GoBack" bookmark marks the author's last edit location; relocate it
# to sit right after the freshly typed "eficacia" (its correct new spot)
# instead of leaving it at the old end-of-paragraph offset.
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()
$bmPos = $d.Range($newWord.End, $newWord.End)
$d.Bookmarks.Add("_GoBack", $bmPos)

# ---------------------------------------------------------------------------
# Changes 2-6: pure run-merge clean-ups (no visible text change).
# ---------------------------------------------------------------------------
Merge-Range "La experiencia realizando el trabajo fue desafiante" "reales de trabajo en equipo."

Merge-Range "Fue una experiencia bastante curiosa" "Fue una buena experiencia como tal"

Merge-Range "La verdad estuvo muy entretenido" "logra el trabajo "

Merge-Range ", med" "fallar y buscar, muy buena "

$ademas = $d.Content
$ademas.Find.Execute("ademas", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$afterAdemas = $ademas.End
$endRng = $d.Content
$endRng.Find.Execute("el trabajo se hace ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$afterHace = $endRng.End
$rng5 = $d.Range($afterAdemas, $afterHace)
$orig5 = $rng5.Text
$rng5.Text = $orig5 + "%"
$rng5b = $d.Range($afterAdemas, $afterAdemas + $orig5.Length + 1)
$rng5b.Text = $orig5

Merge-Range "conflictos su" "conversando y "
